# The underlying species-observation records for rows 2,3,4,5,7,9 on sheet
# "Artfynd" are cyclically rotated among each other:
#   row2 <- old row9, row3 <- old row2, row4 <- old row5,
#   row5 <- old row3, row7 <- old row4, row9 <- old row7
# Rows 1, 6 and 8 are untouched. Columns C,N,P,S,T,U,V,W,Y,Z,AA,AB,AD,AE,AG,
# AT,AW,AX,AY are identical across all these rows, so only the differing
# columns (A,B,D,E,F,G,H,I,J,K,L,Q,R,AF) need to be written.
#
# Column I ("Antal") stores numeric-looking values as text, and K/L/AF are
# sparse text columns where some rows have a present-but-empty cell and
# others have no cell at all; a leading "'" forces Excel to keep a
# numeric-looking or empty value as text instead of coercing it to a blank
#/ number, which is what reproduces that presence distinction.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row 9's data)
$ws.Range("A2").Value = 110282848
$ws.Range("B2").Value = 96348
$ws.Range("D2").Value = "VU"
$ws.Range("E2").Value = 220787
$ws.Range("F2").Value = "Knärot"
$ws.Range("G2").Value = "Goodyera repens"
$ws.Range("H2").Value = "(L.) R. Br."
$ws.Range("I2").Value = "'5"
$ws.Range("J2").Value = "plantor/tuvor"
$ws.Range("K2").Value = "fullt utvecklade blad"
$ws.Range("L2").Value = "'"
$ws.Range("Q2").Value = 600839.9318167433
$ws.Range("R2").Value = 6613983.990819811
$ws.Range("AF2").ClearContents()

# Row 3 (was row 2's data)
$ws.Range("A3").Value = 110282828
$ws.Range("B3").Value = 89425
$ws.Range("D3").Value = "NT"
$ws.Range("E3").Value = 5442
$ws.Range("F3").Value = "Tallticka"
$ws.Range("G3").Value = "Porodaedalea pini"
$ws.Range("H3").Value = "(Brot.) Murrill"
$ws.Range("I3").Value = "'2"
$ws.Range("J3").Value = "fruktkroppar"
$ws.Range("K3").Value = "'"
$ws.Range("L3").ClearContents()
$ws.Range("Q3").Value = 600787.8656294679
$ws.Range("R3").Value = 6613904.709995793
$ws.Range("AF3").Value = "'"

# Row 4 (was row 5's data)
$ws.Range("A4").Value = 110282764
$ws.Range("B4").Value = 96348
$ws.Range("D4").Value = "VU"
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = "Knärot"
$ws.Range("G4").Value = "Goodyera repens"
$ws.Range("H4").Value = "(L.) R. Br."
$ws.Range("I4").Value = "'10"
$ws.Range("J4").Value = "plantor/tuvor"
$ws.Range("K4").Value = "fullt utvecklade blad"
$ws.Range("L4").Value = "'"
$ws.Range("Q4").Value = 600749.0751519018
$ws.Range("R4").Value = 6613971.934424319

# Row 5 (was row 3's data)
$ws.Range("A5").Value = 110282846
$ws.Range("B5").Value = 103288
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 221144
$ws.Range("F5").Value = "Grönpyrola"
$ws.Range("G5").Value = "Pyrola chlorantha"
$ws.Range("H5").Value = "Sw."
$ws.Range("K5").Value = "blomning"
$ws.Range("Q5").Value = 600839.9318167433
$ws.Range("R5").Value = 6613983.990819811
$ws.Range("AF5").ClearContents()

# Row 7 (was row 4's data)
$ws.Range("A7").Value = 110282820
$ws.Range("B7").Value = 89425
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 5442
$ws.Range("F7").Value = "Tallticka"
$ws.Range("G7").Value = "Porodaedalea pini"
$ws.Range("H7").Value = "(Brot.) Murrill"
$ws.Range("Q7").Value = 600724.7123983201
$ws.Range("R7").Value = 6614086.574870056
$ws.Range("AF7").Value = "'"

# Row 9 (was row 7's data)
$ws.Range("A9").Value = 110282856
$ws.Range("B9").Value = 89802
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 5420
$ws.Range("F9").Value = "Grovticka"
$ws.Range("G9").Value = "Phaeolus schweinitzii"
$ws.Range("H9").Value = "(Fr.) Pat."
$ws.Range("I9").Value = "'1"
$ws.Range("J9").Value = "fruktkroppar"
$ws.Range("K9").Value = "'"
$ws.Range("L9").ClearContents()
$ws.Range("Q9").Value = 600677.6983460309
$ws.Range("R9").Value = 6613951.301940188
